$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column H (adds "NRMSE" to the shared strings table)
$ws.Range("H1").Value = "NRMSE"

# Updated data values (rows 2-6, columns A-H)
$data = @(
    @(4.8250205193736733, 4.2065677655677662, 0.49009782689668291, 0.70006987286747524, 0.72427130949858232, 0.50990217310331709, 0.81455566098744159, 0.30654514100213937),
    @(4.8315752399344785, 3.8697759259259259, 1.0520297481319136, 1.0256850140915161, 1.0358072606868107, -0.052029748131913633, 0.0016253820736346842, 0.36110427802200884),
    @(3.4651575951410609, 3.037849276583791, 0.71143510683768885, 0.84346612666881227, 0.87697727384058621, 0.28856489316231115, 0.72563142866077868, 0.28924520827554767),
    @(3.5724667420729541, 3.0515076264880938, 0.80604023236770828, 0.89779743392800371, 0.97058130613488969, 0.19395976763229172, 0.84687427588894859, 0.24335604510033751),
    @(3.3285098047875632, 2.6083217862252175, 0.34801904608118711, 0.5899313909949081, 0.56801432626855775, 0.65198095391881283, 0.87289665666725036, 0.19742051036699659)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rowData = $data[$i]
    for ($j = 0; $j -lt $rowData.Length; $j++) {
        $col = $j + 1
        $ws.Cells.Item($row, $col).Value = $rowData[$j]
    }
}

# Column widths (F and G grew wider to fit the new values; H gets the
# standard column width used elsewhere on the sheet)
$ws.Columns.Item(6).ColumnWidth = 13.666666666666666
$ws.Columns.Item(7).ColumnWidth = 13.833333333333334
$ws.Columns.Item(8).ColumnWidth = 11.833333333333334
